# Add more rows of data (rows 88-98) to Sheet1, and set the sheet view zoom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(323, 1, 5, 1, 6, 7, 6, 7, 4, 5, 1),
    @(253, 1, 7, 6, 5, 4, 4, 6, 3, 6, 4),
    @(85,  7, 6, 1, 7, 1, 2, 7, 7, 5, 2),
    @(153, 2, 1, 1, 5, 1, 2, 2, 1, 3, 1),
    @(139, 6, 1, 1, 5, 1, 6, 7, 6, 2, 6),
    @(246, 7, 4, 1, 5, 2, 3, 5, 5, 3, 4),
    @(156, 3, 4, 4, 6, 1, 7, 7, 1, 2, 3),
    @(237, 7, 4, 2, 6, 2, 7, 4, 7, 6, 6),
    @(274, 1, 4, 4, 6, 4, 4, 7, 4, 6, 1),
    @(308, 7, 4, 4, 6, 6, 7, 4, 5, 2, 4),
    @(200, 7, 6, 3, 5, 7, 7, 6, 2, 2, 3)
)

$startRow = 88
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowVals = $data[$i]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $rowVals[$c]
    }
}

# Set the zoom scale of the active window/view to 85%, as recorded in the diff.
$excel.ActiveWindow.Zoom = 85
